$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 191; this shifts rows 191:264 down to 192:265
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new record's data.
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 44636
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100112017
$ws.Range("G191").Value = "Apio"
$ws.Range("H191").Value = "Americana (o)"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 80
$ws.Range("K191").Value = 12000
$ws.Range("L191").Value = 12000
$ws.Range("M191").Value = 12000
$ws.Range("N191").Value = "`$/docena de matas"
$ws.Range("O191").Value = "Provincia del Elquí"
$ws.Range("P191").Value = 2000
$ws.Range("Q191").Value = 6
$ws.Range("R191").Value = "Hortaliza"
